# Working version of branch and price with L-shaped subproblems
$wb = $excel.ActiveWorkbook

# --- Sheet "deploy_amounts" (row 2 values) ---
$wsAmounts = $wb.Worksheets.Item("deploy_amounts")

$wsAmounts.Range("A2").Value = 80000
$wsAmounts.Range("T2").Value = 0
$wsAmounts.Range("AD2").Value = 80000
$wsAmounts.Range("AP2").Value = 0
$wsAmounts.Range("BC2").Value = 80000

# --- Sheet "deploy_bins" (binary indicator column) ---
$wsBins = $wb.Worksheets.Item("deploy_bins")

$wsBins.Range("A21").Value = 0
$wsBins.Range("A31").Value = 1
$wsBins.Range("A43").Value = 0
$wsBins.Range("A56").Value = 1
